# Insert a new weekly record at row 7 ("Fruta / hortaliza, semanal"):
# this pushes the existing rows 7..26 down to 8..27 and fills the
# newly opened row 7 with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 7, shifting rows 7-26 -> 8-27.
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the new weekly data point.
$ws.Cells.Item(7, 1).Value = 7
$ws.Cells.Item(7, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(7, 3).Value = "Ñuble"
$ws.Cells.Item(7, 4).Value = 44620
$ws.Cells.Item(7, 5).Value = 16
$ws.Cells.Item(7, 6).Value = 100112040
$ws.Cells.Item(7, 7).Value = "Cilantro"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 160
$ws.Cells.Item(7, 11).Value = 550
$ws.Cells.Item(7, 12).Value = 600
$ws.Cells.Item(7, 13).Value = 575
$ws.Cells.Item(7, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(7, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(7, 16).Value = 575
$ws.Cells.Item(7, 17).Value = 1
$ws.Cells.Item(7, 18).Value = "Hortaliza"
